$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: Thursday session -> Promises / promise, promise chaning, fetch
# Row 34: Friday session    -> jQuery   / jquery introduction
$ws.Range("E33").Value = "promise, promise chaning, fetch"
$ws.Range("E34").Value = "jquery introduction"
$ws.Range("D34").Value = "jQuery"
$ws.Range("D33").Value = "Promises"

# Restore the colour-coded "Class Status" / "Session Heading" fills that
# the rest of the sheet uses for completed sessions.
$ws.Range("A33").Interior.Color = 5287936   # green (FF00B050) - matches A26/A28/A30/A31
$ws.Range("A34").Interior.Color = 15773696  # blue  (FF00B0F0) - matches A27/A32

$ws.Range("D33").Interior.Color = 65535     # yellow (FFFFFF00) - matches D26/D28/D30/D32
$ws.Range("D34").Interior.Color = 65535     # yellow (FFFFFF00) - matches D27

# Selection moves from E32 to E34 (author was editing row 34 last)
$ws.Range("E34").Select()
